# Apply the Alvearie FHIR IG metadata refresh to
# StructureDefinition-employee-company.xlsx
#
# Sheet "Metadata" (sheet1):
#   - Version 5.0.0 -> 6.0.0
#   - Date refreshed
#   - Publisher now filled in ("Alvearie Team")
#   - the duplicated "Contact / No display for ContactDetail" row is
#     replaced by a single "Jurisdiction / United States of America" row
#     (the old duplicate row is removed, shrinking the table from 21 to
#     20 data+header rows)
#
# Sheet "Elements" (sheet2):
#   - the root Extension row's Short/Definition text is refreshed

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------

# Remove the second, duplicate "Contact" row (row 11) so the remaining
# rows shift up by one, matching the new A1:B20 extent.
$ws1.Rows.Item(11).Delete()

# Version
$ws1.Range("B3").Value = "6.0.0"

# Date
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher (was blank)
$ws1.Range("B9").Value = "Alvearie Team"

# What used to be the first "Contact" row becomes "Jurisdiction"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Elements sheet ---------------------------------------------------

# Root Extension row ("Short" / "Definition" columns)
$ws2.Range("K2").Value = "Employee Company"
$ws2.Range("L2").Value = "Code for the company of the employee"

Write-Host "edit complete"
